# Change the highlight colour of the "If Landlord or Agent" bullet group
# (the three list items describing the Landlord/Agent dashboard view)
# from light gray to cyan (wdTurquoise), matching both the paragraph mark
# run properties and the text run properties.

$d = $word.ActiveDocument

# Exact paragraph texts (without the trailing paragraph mark) that must be
# re-highlighted. These uniquely identify the three bullets touched by the
# edit, without disturbing the following "View Prospects Notifications"
# bullet (which keeps its light gray highlight) or any other paragraph in
# the document that happens to share similar wording.
$targets = @(
    "If Landlord or Agent: ",
    "View own Property Listings (if available)",
    "Add Property (Else)"
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    $trimmed = $text.TrimEnd([char]13, [char]7)

    foreach ($target in $targets) {
        if ($trimmed -eq $target) {
            # Setting Font.HighlightColorIndex on the paragraph's Range
            # updates the highlight on every run in the paragraph as well
            # as the paragraph mark's run properties (w:pPr/w:rPr), exactly
            # mirroring the lightGray -> cyan change in the diff.
            $p.Range.Font.HighlightColorIndex = 3   # wdTurquoise (cyan)
        }
    }
}
